$wb = $excel.ActiveWorkbook

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4411.185
$ws.Range("I74").Value = 3912.1667
$ws.Range("J74").Value = 4810.4
$ws.Range("K74").Value = 3912.1667
$ws.Range("L74").Value = 4810.4
$ws.Range("M74").Value = -2976.1667
$ws.Range("N74").Value = -6682.4

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4411.185
$ws.Range("I77").Value = 3912.1667
$ws.Range("J77").Value = 4810.4
$ws.Range("K77").Value = 19560.8335
$ws.Range("L77").Value = 24052
$ws.Range("M77").Value = -14880.8335
$ws.Range("N77").Value = -33412

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3369.1738
$ws.Range("I113").Value = 2140.4167
$ws.Range("J113").Value = 4709.636
$ws.Range("K113").Value = 2140.4167
$ws.Range("L113").Value = 4709.636
$ws.Range("M113").Value = 1113.5833
$ws.Range("N113").Value = -11217.636

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3455.8572
$ws.Range("I116").Value = 3894.0952
$ws.Range("K116").Value = 3894.0952
$ws.Range("M116").Value = -452.0952000000002

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 53165.156
$ws.Range("I118").Value = 67122.60000000001
$ws.Range("J118").Value = 824.75
$ws.Range("K118").Value = 201367.8
$ws.Range("L118").Value = 2474.25
$ws.Range("M118").Value = -199710.8
$ws.Range("N118").Value = -5788.25

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2612.0615
$ws.Range("I132").Value = 1377.8043
$ws.Range("J132").Value = 5600.263
$ws.Range("K132").Value = 4133.4129
$ws.Range("L132").Value = 16800.789
$ws.Range("M132").Value = -1603.4129
$ws.Range("N132").Value = -21860.789

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2257.566
$ws.Range("I137").Value = 2219.5642
$ws.Range("J137").Value = 2363.4285
$ws.Range("K137").Value = 6658.692599999999
$ws.Range("L137").Value = 7090.2855
$ws.Range("M137").Value = -4108.692599999999
$ws.Range("N137").Value = -12190.2855

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1499.6285
$ws.Range("I138").Value = 950.1842
$ws.Range("J138").Value = 2152.0938
$ws.Range("K138").Value = 2850.5526
$ws.Range("L138").Value = 6456.2814
$ws.Range("M138").Value = 2289.4474
$ws.Range("N138").Value = -16736.2814

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 669.325
$ws.Range("I141").Value = 669.325
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2007.975
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3172.025
$ws.Range("N141").ClearContents()

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 330345.53
$ws.Range("I61").Value = 257347.39
$ws.Range("J61").Value = 479817.9
$ws.Range("K61").Value = 257347.39
$ws.Range("L61").Value = 479817.9
$ws.Range("M61").Value = -257135.39
$ws.Range("N61").Value = -480241.9

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 139159.84
$ws.Range("I74").Value = 157102.4
$ws.Range("J74").Value = 67389.56
$ws.Range("K74").Value = 157102.4
$ws.Range("L74").Value = 67389.56
$ws.Range("M74").Value = -156228.4
$ws.Range("N74").Value = -69137.56

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 139159.84
$ws.Range("I77").Value = 157102.4
$ws.Range("J77").Value = 67389.56
$ws.Range("K77").Value = 785512
$ws.Range("L77").Value = 336947.8
$ws.Range("M77").Value = -781144
$ws.Range("N77").Value = -345683.8

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2325.7917
$ws.Range("I132").Value = 2075.6182
$ws.Range("J132").Value = 3135.1765
$ws.Range("K132").Value = 6226.8546
$ws.Range("L132").Value = 9405.529500000001
$ws.Range("M132").Value = -3696.8546
$ws.Range("N132").Value = -14465.5295

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 29025.857
$ws.Range("J133").Value = 29025.857
$ws.Range("L133").Value = 29025.857
$ws.Range("N133").Value = -34085.857

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 18833.363
$ws.Range("J135").Value = 18833.363
$ws.Range("L135").Value = 18833.363
$ws.Range("N135").Value = -28973.363

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 330345.53
$ws.Range("I136").Value = 257347.39
$ws.Range("J136").Value = 479817.9
$ws.Range("K136").Value = 772042.17
$ws.Range("L136").Value = 1439453.7
$ws.Range("M136").Value = -769492.17
$ws.Range("N136").Value = -1444553.7

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1630.1904
$ws.Range("I94").Value = 1027.9166
$ws.Range("J94").Value = 2433.2222
$ws.Range("K94").Value = 1027.9166
$ws.Range("L94").Value = 2433.2222
$ws.Range("M94").Value = -576.9166
$ws.Range("N94").Value = -3335.2222

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2295.8088
$ws.Range("I134").Value = 2105.3914
$ws.Range("J134").Value = 2693.9546
$ws.Range("K134").Value = 6316.174199999999
$ws.Range("L134").Value = 8081.8638
$ws.Range("M134").Value = -3781.174199999999
$ws.Range("N134").Value = -13151.8638

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 753.3333
$ws.Range("I16").Value = 725.7143
$ws.Range("J16").Value = 850
$ws.Range("K16").Value = 725.7143
$ws.Range("L16").Value = 850
$ws.Range("M16").Value = -438.7143
$ws.Range("N16").Value = -1424

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3046.7754
$ws.Range("I31").Value = 2150.4517
$ws.Range("J31").Value = 4590.4443
$ws.Range("K31").Value = 2150.4517
$ws.Range("L31").Value = 4590.4443
$ws.Range("M31").Value = -1855.4517
$ws.Range("N31").Value = -5180.4443

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3046.7754
$ws.Range("I34").Value = 2150.4517
$ws.Range("J34").Value = 4590.4443
$ws.Range("K34").Value = 2150.4517
$ws.Range("L34").Value = 4590.4443
$ws.Range("M34").Value = -1948.4517
$ws.Range("N34").Value = -4994.4443

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2672.4922
$ws.Range("I58").Value = 2661.5
$ws.Range("J58").Value = 2716.4614
$ws.Range("K58").Value = 2661.5
$ws.Range("L58").Value = 2716.4614
$ws.Range("M58").Value = -2458.5
$ws.Range("N58").Value = -3122.4614

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 753.3333
$ws.Range("I113").Value = 725.7143
$ws.Range("J113").Value = 850
$ws.Range("K113").Value = 725.7143
$ws.Range("L113").Value = 850
$ws.Range("M113").Value = 1444.2857
$ws.Range("N113").Value = -5190

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2053.6072
$ws.Range("I132").Value = 1075.05
$ws.Range("K132").Value = 3225.15
$ws.Range("M132").Value = -695.1499999999996

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2672.4922
$ws.Range("I136").Value = 2661.5
$ws.Range("J136").Value = 2716.4614
$ws.Range("K136").Value = 7984.5
$ws.Range("L136").Value = 8149.3842
$ws.Range("M136").Value = -5434.5
$ws.Range("N136").Value = -13249.3842

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 470.675
$ws.Range("I5").Value = 404.74194
$ws.Range("J5").Value = 697.7778
$ws.Range("K5").Value = 1214.22582
$ws.Range("L5").Value = 2093.3334
$ws.Range("M5").Value = -1102.22582
$ws.Range("N5").Value = -2317.3334

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 169005.5
$ws.Range("J130").Value = 169005.5
$ws.Range("L130").Value = 507016.5
$ws.Range("N130").Value = -517056.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1202.2097
$ws.Range("I131").Value = 1259.8334
$ws.Range("J131").Value = 1188.38
$ws.Range("K131").Value = 3779.5002
$ws.Range("L131").Value = 3565.14
$ws.Range("M131").Value = 1260.4998
$ws.Range("N131").Value = -13645.14

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 470.675
$ws.Range("I135").Value = 404.74194
$ws.Range("J135").Value = 697.7778
$ws.Range("K135").Value = 3642.67746
$ws.Range("L135").Value = 6280.000199999999
$ws.Range("M135").Value = -1107.67746
$ws.Range("N135").Value = -11350.0002

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -4540

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -5872

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 5548.8423
$ws.Range("I107").Value = 10149.4
$ws.Range("J107").Value = 437.1111
$ws.Range("K107").Value = 10149.4
$ws.Range("L107").Value = 437.1111
$ws.Range("M107").Value = -8229.4
$ws.Range("N107").Value = -4277.1111

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1503.2162
$ws.Range("I93").Value = 1448.76
$ws.Range("J93").Value = 1616.6666
$ws.Range("K93").Value = 1448.76
$ws.Range("L93").Value = 1616.6666
$ws.Range("M93").Value = -200.76
$ws.Range("N93").Value = -4112.6666

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6653.143
$ws.Range("I132").Value = 2305.8
$ws.Range("J132").Value = 13046.294
$ws.Range("K132").Value = 6917.400000000001
$ws.Range("L132").Value = 39138.882
$ws.Range("M132").Value = -4387.400000000001
$ws.Range("N132").Value = -44198.882

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 47330
$ws.Range("J140").Value = 47330
$ws.Range("L140").Value = 47330
$ws.Range("N140").Value = -57690

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 467.875
$ws.Range("I113").Value = 478.6
$ws.Range("J113").Value = 450
$ws.Range("K113").Value = 1435.8
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 734.1999999999998
$ws.Range("N113").Value = -5690

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1870
$ws.Range("I132").Value = 1265
$ws.Range("J132").Value = 3080
$ws.Range("K132").Value = 3795
$ws.Range("L132").Value = 9240
$ws.Range("M132").Value = -1265
$ws.Range("N132").Value = -14300
